# MAD Lab2 and CSAD Lab2 created
#
# On the "significant figures" worked-example slide, the third example
# ("8.76 rounded to 8.82") had its "s.f." note merged onto the same line
# with the figure count ("2") added in front of it.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(29)
$shp = $s.Shapes.Item(2)
$tf = $shp.TextFrame
$tr = $tf.TextRange

$oldText = "`t   8.76 rounded to 8.82 s.f."
$newText = "`t   8.76 rounded to 8.82    2s.f."

[void]$tr.Replace($oldText, $newText)
